# Generate Report for Handback
# Replace the two UUID-named e2e/*.md entries that this CI run tracked with
# the two new ones from this run, refresh the recorded timestamps and the
# combined handoff xliff filenames, and collapse the now-merged per-file
# xliff rows down to the single generated xliff per locale.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "2a8a0528-1b74-4041-a344-c4b731c644b3"
$oldUuid2 = "d777324f-c81b-4580-bd8e-123f1e752530"
$newUuid1 = "b379717d-78e0-454a-966f-ac667ff62196"
$newUuid2 = "ffff964063b8-d866-4ebf-af8b-75f6ecefa83c"
$newHash  = "fd3b16db992c7dc4ca507493b5492e41a4af7254"

$file1Name = "$newUuid1.md"
$file1Path = "e2e\$newUuid1.md"
$file2Name = "$newUuid2.md"
$file2Path = "e2e\$newUuid2.md"

$genDate        = "2016-08-25 23:00:52"
$zhcnXlfName    = "$newUuid1.$newHash.zh-cn.xlf"
$zhcnHandoffDt  = "2016-08-25 23:00:48"
$zhcnHandbackDt = "2016-08-25 23:01:21"
$dedeXlfName    = "$newUuid1.$newHash.de-de.xlf"
$dedeHandbackDt = "2016-08-25 23:01:28"

$hyperlinkBase1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38a36c9885fb1af62d0228f32eb5ba585fcbeae8/e2e/$oldUuid1.md"
$hyperlinkBase2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38a36c9885fb1af62d0228f32eb5ba585fcbeae8/e2e/$oldUuid2.md"
$hyperlinkZhcn1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d73012c6fe273e944d3461a7214b40fa67b041a5/e2e/$oldUuid1.md"
$hyperlinkZhcn2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d73012c6fe273e944d3461a7214b40fa67b041a5/e2e/$oldUuid2.md"
$hyperlinkDede1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/07875ce156ba18f4bad3a3d33b9a03b69b1e1653/e2e/$oldUuid1.md"
$hyperlinkDede2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/07875ce156ba18f4bad3a3d33b9a03b69b1e1653/e2e/$oldUuid2.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name / Path And Name(hyperlink) / Extension /
# Publish URL / zh-cn / de-de / Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $file1Name
$wsOverview.Range("B2").Value = $file1Path
$wsOverview.Range("G2").Value = $genDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkBase1, "", "", $file1Path)
$wsOverview.Range("B2").Style = "HyperLink"

$wsOverview.Range("A3").Value = $file2Name
$wsOverview.Range("B3").Value = $file2Path
$wsOverview.Range("G3").Value = $genDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkBase2, "", "", $file2Path)
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name(A,I hyperlinks) / Correspond Handoff
# File(G,J) / Correspond Handoff Datetime(H) / Correspond Handback
# DateTime(K)
# ---------------------------------------------------------------------------
$wsZhcn = $wb.Worksheets.Item("zh-cn")
$wsZhcn.Hyperlinks.Delete()

$wsZhcn.Range("A2").Value = $file1Name
$wsZhcn.Range("I2").Value = $file1Name
$wsZhcn.Range("G2").Value = $zhcnXlfName
$wsZhcn.Range("J2").Value = $zhcnXlfName
$wsZhcn.Range("H2").Value = $zhcnHandoffDt
$wsZhcn.Range("K2").Value = $zhcnHandbackDt
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("A2"), $hyperlinkBase1, "", "", $file1Name)
$wsZhcn.Range("A2").Style = "HyperLink"
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("I2"), $hyperlinkZhcn1, "", "", $file1Name)
$wsZhcn.Range("I2").Style = "HyperLink"

$wsZhcn.Range("A3").Value = $file2Name
$wsZhcn.Range("I3").Value = $file2Name
$wsZhcn.Range("G3").Value = $zhcnXlfName
$wsZhcn.Range("J3").Value = $zhcnXlfName
$wsZhcn.Range("H3").Value = $zhcnHandoffDt
$wsZhcn.Range("K3").Value = $zhcnHandbackDt
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("A3"), $hyperlinkBase2, "", "", $file2Name)
$wsZhcn.Range("A3").Style = "HyperLink"
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("I3"), $hyperlinkZhcn2, "", "", $file2Name)
$wsZhcn.Range("I3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "de-de": Source File Name(A,I hyperlinks) / Correspond Handoff
# File(G,J) / Correspond Handoff Datetime(H) / Correspond Handback
# DateTime(K)
# ---------------------------------------------------------------------------
$wsDede = $wb.Worksheets.Item("de-de")
$wsDede.Hyperlinks.Delete()

$wsDede.Range("A2").Value = $file1Name
$wsDede.Range("I2").Value = $file1Name
$wsDede.Range("G2").Value = $dedeXlfName
$wsDede.Range("J2").Value = $dedeXlfName
$wsDede.Range("H2").Value = $genDate
$wsDede.Range("K2").Value = $dedeHandbackDt
$wsDede.Hyperlinks.Add($wsDede.Range("A2"), $hyperlinkBase1, "", "", $file1Name)
$wsDede.Range("A2").Style = "HyperLink"
$wsDede.Hyperlinks.Add($wsDede.Range("I2"), $hyperlinkDede1, "", "", $file1Name)
$wsDede.Range("I2").Style = "HyperLink"

$wsDede.Range("A3").Value = $file2Name
$wsDede.Range("I3").Value = $file2Name
$wsDede.Range("G3").Value = $dedeXlfName
$wsDede.Range("J3").Value = $dedeXlfName
$wsDede.Range("H3").Value = $genDate
$wsDede.Range("K3").Value = $dedeHandbackDt
$wsDede.Hyperlinks.Add($wsDede.Range("A3"), $hyperlinkBase2, "", "", $file2Name)
$wsDede.Range("A3").Style = "HyperLink"
$wsDede.Hyperlinks.Add($wsDede.Range("I3"), $hyperlinkDede2, "", "", $file2Name)
$wsDede.Range("I3").Style = "HyperLink"
